# Update excel writer with format function
#
# This script updates the "test.xlsx"-style workbook:
#  - Sheet "Version": refresh Start Time / End Time / Run Time strings
#  - Sheet "3_3" (averages): update the "Deep Sleep" row's computed stats and
#    add a new "BT Idle" row with its own computed stats
#  - Sheet "1_8" (offsets): update the "Deep Sleep" row's computed stats and
#    add a new "BT Idle" row with its own computed stats

$wb = $excel.ActiveWorkbook

$wsVersion = $wb.Worksheets.Item("Version")
$wsAvg     = $wb.Worksheets.Item("3_3")
$wsOffset  = $wb.Worksheets.Item("1_8")

# ---------------------------------------------------------------------------
# 1. Refresh the timing information on the "Version" sheet
# ---------------------------------------------------------------------------
$wsVersion.Range("C10").Value = "09:54:10 03-15-2017 Wednesday Pacific Daylight Time"
$wsVersion.Range("C11").Value = "09:54:45 03-15-2017 Wednesday Pacific Daylight Time"
$wsVersion.Range("C12").Value = "0:0:34 0 days"

# ---------------------------------------------------------------------------
# Helper raw-data series (100 samples each) used to compute
# Average / Max / Min / Sdev(population) for each metric/state pair
# ---------------------------------------------------------------------------
$deepSleepAvgRaw = "0.574475,0.602885,0.567342,0.577349,0.576153,0.567892,0.56919,0.593048,0.571506,0.571886,0.565841,0.602512,0.572348,0.566771,0.574475,0.574108,0.595644,0.561302,0.566676,0.575684,0.578755,0.568911,0.571785,0.576798,0.572437,0.567613,0.600012,0.575494,0.571513,0.565657,0.572253,0.578008,0.586078,0.567151,0.574298,0.569196,0.568082,0.568347,0.602518,0.577172,0.575032,0.571411,0.605494,0.605487,0.567056,0.57188,0.574101,0.575868,0.592111,0.570019,0.561941,0.572158,0.581438,0.571791,0.596479,0.594991,0.572342,0.568442,0.580888,0.569101,0.593795,0.562498,0.568911,0.574672,0.566676,0.566683,0.596581,0.570019,0.572715,0.569006,0.613837,0.566778,0.564455,0.567233,0.601771,0.567063,0.559434,0.567518,0.599733,0.602138,0.567335,0.562777,0.577084,0.574577,0.559631,0.573931,0.575868,0.571785,0.566961,0.569835,0.569652,0.565575,0.592023,0.601411,0.573273,0.572899,0.562783,0.575589,0.576513,0.566683"

$btIdleAvgRaw = "9.840451,9.847224,9.848256,9.831905,9.82077,9.830329,9.83664,9.830892,9.810831,9.816497,9.813807,9.825125,9.812231,9.807129,9.797557,9.818447,9.813345,9.801266,9.785675,9.800627,9.81437,9.809439,9.783814,9.777326,9.797006,9.803222,9.799873,9.778813,9.781864,9.800716,9.793664,9.787815,9.773711,9.78289,9.79969,9.788555,9.772869,9.778617,9.792917,9.798202,9.784282,9.768141,9.790308,9.787536,9.789479,9.774812,9.780661,9.784758,9.797278,9.791341,9.772951,9.779547,9.788936,9.796361,9.780383,9.768691,9.777971,9.803229,9.79524,9.77825,9.774901,9.786415,9.798481,9.784561,9.772679,9.785125,9.803786,9.7945,9.781504,9.763868,9.787802,9.79524,9.78615,9.774072,9.777788,9.787346,9.790689,9.78734,9.77863,9.779364,9.788555,9.790593,9.783358,9.782618,9.787163,9.794214,9.787068,9.779642,9.78308,9.787713,9.790417,9.785865,9.771008,9.77806,9.794873,9.792733,9.787074,9.77435,9.787435,9.797842"

$deepSleepOffsetRaw = "0.00393,-0.006885,-0.000238,0.003918,-0.007157,-0.003537,-0.002162,0.000878,0.000311,-0.002205,-0.003278,-0.009956,-0.000793,-0.009679,-0.009679,-0.000503,-0.000793,0.002537,-0.002446,0.001162,-0.003019,0.004454,-0.004141,-0.003586,0.004756,0.003394,-0.003556,-0.007459,-0.004092,0.001137,0.001451,-0.006602,0.002814,0.000348,-0.009148,-0.00247,-0.003309,-0.002156,-0.004129,-0.003833,-0.001619,-0.002723,-0.002446,-0.003007,-0.001083,-0.001909,-0.007175,0.000576,5.2e-05,0.003647,-0.000818,-0.006084,-0.001903,-0.007169,-0.000238,0.001439,-0.003839,-0.002452,-0.002181,-0.004135,0.003375,0.000613,-0.006941,-0.005775,-0.001638,-0.003574,-0.002446,-0.003001,0.000576,-0.001626,0.000588,8e-06,-0.00633,0.002247,-0.001897,-0.004394,-0.00773,0.003067,-0.00054,-0.007989,-0.004672,-0.003303,0.005268,-0.006891,-0.00136,-0.005227,0.000588,-0.006071,0.003936,-0.005492,-0.005504,0.003079,-0.000528,-0.007187,-0.000516,3.3e-05,-0.002476,-0.004394,-0.004425,-0.003586"

$btIdleOffsetRaw = "0.001445,-0.001077,-0.008316,-0.00276,-0.005245,0.00755,-0.007163,-0.001354,0.000323,0.002543,-0.003297,-0.003278,-0.000251,-0.010782,0.003369,-0.008556,0.002845,-0.002187,0.003079,-0.002464,0.000619,-0.004961,-0.003278,0.000317,-0.001903,-0.001922,-0.003013,-0.003858,-0.004702,0.001717,0.003388,0.000335,3.3e-05,-0.004141,0.00364,-0.001354,-0.003031,-0.007749,-0.004111,-0.001903,-0.001909,0.001988,-0.005516,-0.00329,-0.006367,-0.000775,-0.002742,-0.007194,-0.003556,-0.003827,-0.003272,0.00171,0.000859,-0.004129,-0.003253,-0.006614,-0.006374,0.000859,-0.001595,-0.007718,0.001427,-0.001891,-0.012188,-0.001909,-0.006904,-0.0072,0.002537,0.000582,-0.004659,0.000878,0.004473,-0.000214,4.5e-05,-0.007453,-0.004968,-0.005523,-0.000522,-0.007995,-0.007477,-0.001052,-0.00329,-0.00244,-0.006337,-0.007459,-0.002729,-0.002748,-0.004974,-0.003852,-0.000799,-0.0044,-0.005763,-0.002988,-0.006898,-0.003858,-0.009389,-0.00358,0.000582,0.000613,-0.007471,-0.001348"

# ---------------------------------------------------------------------------
# Helper that, given a comma separated list of numbers, returns a hashtable
# with the rounded (3 decimal) Average / Max / Min / population-Sdev stats,
# mirroring the "format function" used by the excel writer.
# ---------------------------------------------------------------------------
function Get-Stats([string]$raw) {
    $nums = $raw.Split(",") | ForEach-Object { [double]$_ }
    $count = $nums.Count

    $sum = 0
    foreach ($n in $nums) { $sum += $n }
    $avg = $sum / $count

    $max = $nums[0]
    $min = $nums[0]
    foreach ($n in $nums) {
        if ($n -gt $max) { $max = $n }
        if ($n -lt $min) { $min = $n }
    }

    $sq = 0
    foreach ($n in $nums) { $sq += [math]::Pow(($n - $avg), 2) }
    $sdev = [math]::Sqrt($sq / $count)

    return @{
        Average = [math]::Round($avg, 3)
        Max     = [math]::Round($max, 3)
        Min     = [math]::Round($min, 3)
        Sdev    = [math]::Round($sdev, 3)
        Count   = $count
    }
}

$deepSleepAvgStats    = Get-Stats $deepSleepAvgRaw
$btIdleAvgStats       = Get-Stats $btIdleAvgRaw
$deepSleepOffsetStats = Get-Stats $deepSleepOffsetRaw
$btIdleOffsetStats    = Get-Stats $btIdleOffsetRaw

# ---------------------------------------------------------------------------
# 2. Sheet "3_3" -- averages
# ---------------------------------------------------------------------------

# Refresh the existing "Deep Sleep" row (row 3)
$wsAvg.Range("C3").Value = $deepSleepAvgStats.Average
$wsAvg.Range("D3").Value = $deepSleepAvgStats.Max
$wsAvg.Range("E3").Value = $deepSleepAvgStats.Min
$wsAvg.Range("F3").Value = $deepSleepAvgStats.Sdev
$wsAvg.Range("G3").Value = $deepSleepAvgStats.Count
$wsAvg.Range("H3").Value = $deepSleepAvgRaw

# Add the new "BT Idle" row (row 4), cloning the formatting of row 3
$wsAvg.Range("B3:H3").Copy($wsAvg.Range("B4:H4"))
$wsAvg.Range("B4").Value = "BT Idle"
$wsAvg.Range("C4").Value = $btIdleAvgStats.Average
$wsAvg.Range("D4").Value = $btIdleAvgStats.Max
$wsAvg.Range("E4").Value = $btIdleAvgStats.Min
$wsAvg.Range("F4").Value = $btIdleAvgStats.Sdev
$wsAvg.Range("G4").Value = $btIdleAvgStats.Count
$wsAvg.Range("H4").Value = $btIdleAvgRaw

# ---------------------------------------------------------------------------
# 3. Sheet "1_8" -- offsets
# ---------------------------------------------------------------------------

# Refresh the existing "Deep Sleep" row (row 3)
$wsOffset.Range("C3").Value = $deepSleepOffsetStats.Average
$wsOffset.Range("D3").Value = $deepSleepOffsetStats.Max
$wsOffset.Range("E3").Value = $deepSleepOffsetStats.Min
$wsOffset.Range("F3").Value = $deepSleepOffsetStats.Sdev
$wsOffset.Range("G3").Value = $deepSleepOffsetStats.Count
$wsOffset.Range("H3").Value = $deepSleepOffsetRaw

# Add the new "BT Idle" row (row 4), cloning the formatting of row 3
$wsOffset.Range("B3:H3").Copy($wsOffset.Range("B4:H4"))
$wsOffset.Range("B4").Value = "BT Idle"
$wsOffset.Range("C4").Value = $btIdleOffsetStats.Average
$wsOffset.Range("D4").Value = $btIdleOffsetStats.Max
$wsOffset.Range("E4").Value = $btIdleOffsetStats.Min
$wsOffset.Range("F4").Value = $btIdleOffsetStats.Sdev
$wsOffset.Range("G4").Value = $btIdleOffsetStats.Count
$wsOffset.Range("H4").Value = $btIdleOffsetRaw
